$d = $word.ActiveDocument

# Paragraph 1 is the title: "Retail Sales Data Forecasting Models"
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

# Bump the title's font size from 32 half-points (16pt) to 40 half-points (20pt),
# updating both the regular and complex-script (Cs) sizes.
$titleRange.Font.Size = 20
$titleRange.Font.SizeBi = 20

# Insert a new paragraph containing "Dubstech Datathon 2020" right after the
# title and before the "Alexander Van Roijen, Ashley Batchelor" paragraph.
# Inserting before paragraph 2 makes the new paragraph inherit that
# paragraph's (non-bold, 32 half-point) formatting.
$secondPara = $d.Paragraphs.Item(2)
$insertPoint = $secondPara.Range
$insertPoint.Collapse(1)  # wdCollapseStart
$insertPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range
$newRange.Text = "Dubstech Datathon 2020"
